$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.887.73"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "1.889.94"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'0.7694"
$ws.Range("E5").Value = "  -0.66%  "
$ws.Range("D6").Value = "'242.89"
$ws.Range("E6").Value = "  -0.77%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "'0.3128"
$ws.Range("E8").Value = "  -0.33%  "
$ws.Range("D9").Value = "'25.73"
$ws.Range("E9").Value = "  +1.61%  "
$ws.Range("D10").Value = "'0.07182"
$ws.Range("E10").Value = "  -3.50%  "
$ws.Range("D11").Value = "'0.08578"
$ws.Range("E11").Value = "  +5.54%  "
$ws.Range("D12").Value = "'0.7652"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").Value = "1.930.73"
$ws.Range("E13").Value = "  +2.59%  "
$ws.Range("D14").Value = "'5.369"
$ws.Range("E14").Value = "  -1.77%  "
$ws.Range("D15").Value = "'93.76"
$ws.Range("E15").Value = "  +1.66%  "
$ws.Range("D16").Value = "'6.163"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").Value = "29.931.16"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("D18").Value = "'13.79"
$ws.Range("E18").Value = "  -1.26%  "
$ws.Range("D19").Value = "'244.94"
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("D20").Value = "'0.000007821"
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("D21").Value = "2.169.54"
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("D22").Value = "'0.9994"
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("D23").Value = "'8.053"
$ws.Range("E23").Value = "  -0.55%  "
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").Value = "'0.1642"
$ws.Range("E25").Value = "  +4.11%  "
$ws.Range("D26").Value = "'9.396"
$ws.Range("E26").Value = "  -0.28%  "
$ws.Range("D27").Value = "'162.92"
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("D28").Value = "'18.76"
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("D29").Value = "'2.037"
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("D30").Value = "'1.465"
$ws.Range("E30").Value = "  +1.75%  "
$ws.Range("D31").Value = "'1.536"
$ws.Range("E31").Value = "  -0.81%  "
$ws.Range("D32").Value = "'4.518"
$ws.Range("E32").Value = "  +0.54%  "
$ws.Range("D33").Value = "'4.102"
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("D34").Value = "'0.05455"
$ws.Range("E34").Value = "  -1.06%  "
$ws.Range("D35").Value = "'1.245"
$ws.Range("E35").Value = "  -0.14%  "
$ws.Range("D36").Value = "'0.7454"
$ws.Range("E36").Value = "  -1.76%  "
$ws.Range("D37").Value = "'1.002"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "'2.700"
$ws.Range("E38").Value = "  +2.05%  "
$ws.Range("D39").Value = "'0.01956"
$ws.Range("E39").Value = "  +1.54%  "
$ws.Range("E40").Value = "  -0.18%  "
$ws.Range("D41").Value = "'0.4475"
$ws.Range("E41").Value = "  +0.43%  "
$ws.Range("D42").Value = "1.108.18"
$ws.Range("E42").Value = "  -4.67%  "
$ws.Range("D43").Value = "'73.30"
$ws.Range("E43").Value = "  -0.74%  "
$ws.Range("D44").Value = "'6.081"
$ws.Range("E44").Value = "  +1.96%  "
$ws.Range("D45").Value = "'0.8530"
$ws.Range("E45").Value = "  +0.80%  "
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("D47").Value = "'102.85"
$ws.Range("E47").Value = "  +0.62%  "
$ws.Range("D48").Value = "'7.679"
$ws.Range("E48").Value = "  +1.94%  "
$ws.Range("D49").Value = "'1.868"
$ws.Range("E49").Value = "  -1.78%  "
$ws.Range("D50").Value = "'3.006"
$ws.Range("E50").Value = "  -2.81%  "
$ws.Range("D51").Value = "2.066.30"
$ws.Range("E51").Value = "  +1.17%  "
